$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Locate the "Senior Computer Scientist – Web Standards" bullet paragraph
# (8th paragraph in the Content Placeholder).
$para = $tr.Paragraphs(8)
$paraStart = $para.Start

# The paragraph currently reads:
#   "Senior Computer Scientist – Web Standards"
# Keep the "Senior Computer Scientist – " prefix (characters 1-28) as-is and
# replace the trailing "Web Standards" (characters 29-41, length 13) with the
# new wording, split across two runs:
#   "Web Platform Innovation " + "and Standards"

$tail = $tr.Characters($paraStart + 28, 13)
$tail.Text = "Web Platform Innovation and Standards"

$lastRun = $tr.Characters($paraStart + 28 + 24, 13)
$lastRun.Text = "and Standards"
